$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Number" header above the row-number column.
$ws.Range("A1").Value = "Number"

# Add "Notes" header in the new trailing column, styled to match the
# other bold headers in row 1, with a thin left+right border.
$ws.Range("AH1").Value = "Notes"
$ws.Range("AH1").Font.Bold = $true
$ws.Range("AH1").Borders.Item(7).LineStyle = 1
$ws.Range("AH1").Borders.Item(10).LineStyle = 1

# AG5 held a broken array formula ("=unchanged" => #NAME?); replace it with
# the plain literal text the formula was trying to produce.
$ws.Range("AG5").Value = "unchanged"
